$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting from H1 (bold, bordered, centered header style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: row -> (I value, J value)
$data = @{
    2  = @(1, 2)
    3  = @(1, 6)
    4  = @(1, 4)
    5  = @(1, 7)
    6  = @(1, 5)
    7  = @(1, 4)
    8  = @(1, 2)
    9  = @(1, 4)
    10 = @(1, 7)
    11 = @(1, 7)
    12 = @(1, 7)
    13 = @(4, 7)
    14 = @(2, 4)
    15 = @(1, 2)
    16 = @(1, 3)
    17 = @(5, 7)
    18 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
